$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Add attendance rows for 10-Oct-2023 (row 44) and 11-Oct-2023 (row 45)
# ---------------------------------------------------------------

# Row 44 - 10th Oct 2023 (serial 45209)
$ws.Range("A44").Value = 45209
$ws.Range("A44").NumberFormat = $ws.Range("A43").NumberFormat

$ws.Range("B44").Value = "ABSENT"
$ws.Range("C44").Value = "PRESENT"
$ws.Range("D44").Value = "PRESENT"
$ws.Range("E44").Value = "PRESENT"
$ws.Range("F44").Value = "PRESENT"
$ws.Range("G44").Value = "ABSENT"
$ws.Range("H44").Value = "PRESENT"
$ws.Range("I44").Value = "ABSENT"
$ws.Range("J44").Value = "ABSENT"
$ws.Range("K44").Value = "PRESENT"

# Row 45 - 11th Oct 2023 (serial 45210)
$ws.Range("A45").Value = 45210
$ws.Range("A45").NumberFormat = $ws.Range("A43").NumberFormat

$ws.Range("B45").Value = "PRESENT"
$ws.Range("C45").Value = "PRESENT"
$ws.Range("D45").Value = "PRESENT"
$ws.Range("E45").Value = "PRESENT"
$ws.Range("F45").Value = "PRESENT"
$ws.Range("G45").Value = "ABSENT"
$ws.Range("H45").Value = "ABSENT"
$ws.Range("I45").Value = "ABSENT"
$ws.Range("J45").Value = "ABSENT"
$ws.Range("K45").Value = "PRESENT"

# ---------------------------------------------------------------
# Comments explaining the ABSENT marks
# ---------------------------------------------------------------
$null = $ws.Range("B44").AddComment("Dell:`nDue to some work")
$null = $ws.Range("G44").AddComment("Dell:`nNot informed")
$null = $ws.Range("I44").AddComment("Dell:`nSame")
$null = $ws.Range("J44").AddComment("Dell:`nSame")

$null = $ws.Range("G45").AddComment("Dell:`nNot informed")
$null = $ws.Range("H45").AddComment("Dell:`nDue to exam")
$null = $ws.Range("I45").AddComment("Dell:`nSame")
$null = $ws.Range("J45").AddComment("Dell:`nSame")

# ---------------------------------------------------------------
# Extend the date-column validation so it also covers the new rows
# ---------------------------------------------------------------
$oldDateValidation = $ws.Range("A40:A43")
$oldDateValidation.Validation.Delete()
$newDateValidation = $ws.Range("A40:A45")
$newDateValidation.Validation.Add(0, 1, 1, "")

# ---------------------------------------------------------------
# Leave the cursor where the author left it after the last edit
# ---------------------------------------------------------------
$ws.Range("B48").Select()
